# Atualizado por script em 05-11-2023 08:45
# Adds 3 new match rows (104-106) to the Costa Rica Primera Division sheet,
# mirroring the formatting of the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) of the last existing data row (103) down to
# the three new rows (104-106) so the new cells inherit the same number
# formats / borders used throughout the sheet (bold/bordered index column,
# date-time formatted match-date column).
$ws.Range("A103:V103").Copy()
$ws.Range("A104:V106").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 104
$ws.Range("A104").Value = 103
$ws.Range("B104").Value = "costa-rica"
$ws.Range("C104").Value = "primera-division"
$ws.Range("D104").Value = "2023-2024"
$ws.Range("E104").Value = 45234.91666666666
$ws.Range("F104").Value = "Guanacasteca"
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = "AD Santos"
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 1.88
$ws.Range("K104").Value = "31/10/2023 14:42"
$ws.Range("L104").Value = 1.98
$ws.Range("M104").Value = "04/11/2023 21:54"
$ws.Range("N104").Value = 3.38
$ws.Range("O104").Value = "31/10/2023 14:42"
$ws.Range("P104").Value = 3.34
$ws.Range("Q104").Value = "04/11/2023 21:54"
$ws.Range("R104").Value = 4.21
$ws.Range("S104").Value = "31/10/2023 14:42"
$ws.Range("T104").Value = 4.13
$ws.Range("U104").Value = "04/11/2023 21:54"
$ws.Range("V104").Value = "https://www.betexplorer.com/football/costa-rica/primera-division/guanacasteca-santos-de-guapiles/YHvtcjRB/"

# Row 105
$ws.Range("A105").Value = 104
$ws.Range("B105").Value = "costa-rica"
$ws.Range("C105").Value = "primera-division"
$ws.Range("D105").Value = "2023-2024"
$ws.Range("E105").Value = 45235
$ws.Range("F105").Value = "Sporting San Jose"
$ws.Range("G105").Value = 4
$ws.Range("H105").Value = "Zeledon"
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 1.79
$ws.Range("K105").Value = "29/10/2023 01:12"
$ws.Range("L105").Value = 1.66
$ws.Range("M105").Value = "04/11/2023 23:59"
$ws.Range("N105").Value = 3.72
$ws.Range("O105").Value = "29/10/2023 01:12"
$ws.Range("P105").Value = 3.9
$ws.Range("Q105").Value = "04/11/2023 23:59"
$ws.Range("R105").Value = 4.32
$ws.Range("S105").Value = "29/10/2023 01:12"
$ws.Range("T105").Value = 5.23
$ws.Range("U105").Value = "04/11/2023 23:59"
$ws.Range("V105").Value = "https://www.betexplorer.com/football/costa-rica/primera-division/sporting-san-jose-zeledon/K8DbFXtt/"

# Row 106
$ws.Range("A106").Value = 105
$ws.Range("B106").Value = "costa-rica"
$ws.Range("C106").Value = "primera-division"
$ws.Range("D106").Value = "2023-2024"
$ws.Range("E106").Value = 45235.125
$ws.Range("F106").Value = "Saprissa"
$ws.Range("G106").Value = 1
$ws.Range("H106").Value = "Alajuelense"
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 1.74
$ws.Range("K106").Value = "31/10/2023 14:42"
$ws.Range("L106").Value = 2
$ws.Range("M106").Value = "05/11/2023 02:59"
$ws.Range("N106").Value = 3.71
$ws.Range("O106").Value = "31/10/2023 14:42"
$ws.Range("P106").Value = 3.49
$ws.Range("Q106").Value = "05/11/2023 02:59"
$ws.Range("R106").Value = 4.31
$ws.Range("S106").Value = "31/10/2023 14:42"
$ws.Range("T106").Value = 3.86
$ws.Range("U106").Value = "05/11/2023 02:59"
$ws.Range("V106").Value = "https://www.betexplorer.com/football/costa-rica/primera-division/saprissa-alajuelense/UqYleUeO/"
